# Refresh the crypto "symbol list" price/volume(1h) figures (D/E columns,
# rows 2-51) to match the latest scrape, per the GitHub Actions commit.
# Values are stored as text (not numbers) in the sheet, so we force the
# cell to Text format before assigning, then clear that temporary format
# again so we don't leave a stray numeric/percent display format behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "335.61"
    "E2" = "1.94%"
    "D3" = "43.97"
    "E3" = "6.81%"
    "D4" = "5.718"
    "E4" = "1.58%"
    "D5" = "0.08316"
    "E5" = "1.36%"
    "D6" = "8.857"
    "E6" = "1.27%"
    "D7" = "1.962"
    "E7" = "-2.16%"
    "D8" = "2.879"
    "E8" = "-2.89%"
    "D9" = "0.9432"
    "E9" = "2.30%"
    "D10" = "0.1248"
    "E10" = "-2.57%"
    "D11" = "0.1979"
    "E11" = "1.10%"
    "D12" = "0.1058"
    "E12" = "12.76%"
    "D13" = "0.04708"
    "E13" = "22.25%"
    "E14" = "0.91%"
    "D15" = "0.001290"
    "E15" = "-1.53%"
    "D16" = "0.005925"
    "E16" = "-4.08%"
    "D17" = "3.500"
    "E17" = "1.52%"
    "D18" = "4.523"
    "E18" = "0.47%"
    "D20" = "8.770"
    "E20" = "6.64%"
    "D21" = "0.1353"
    "E21" = "-0.85%"
    "D22" = "0.2692"
    "E22" = "12.68%"
    "D23" = "0.04406"
    "E23" = "-0.01%"
    "E24" = "0.51%"
    "D25" = "0.004391"
    "E25" = "1.83%"
    "E26" = "5.02%"
    "D27" = "0.0003994"
    "E27" = "-94.68%"
    "D39" = "0.02811"
    "E39" = "1.51%"
    "D40" = "0.06079"
    "E40" = "11.30%"
    "D41" = "0.007912"
    "E41" = "-0.54%"
    "D42" = "0.1426"
    "E42" = "0.54%"
    "D43" = "0.008978"
    "D44" = "0.002123"
    "E44" = "-2.29%"
    "D45" = "0.01059"
    "E45" = "-7.42%"
    "D46" = "0.00007009"
    "E46" = "3.54%"
    "E47" = "0.03%"
    "D48" = "0.003189"
    "E48" = "-0.13%"
    "D49" = "0.002273"
    "E49" = "-0.28%"
    "D50" = "0.00002103"
    "E50" = "0.03%"
    "D51" = "0.0002003"
    "E51" = "0.03%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}